$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-format style from the last existing data row (A385) down
# through the new rows so the new dates keep the same number format/border.
$ws.Range("A385").Copy()
$ws.Range("A386:A464").PasteSpecial(-4122)

$rowData = @(
    "386|44460|0|1|53.73455131649651",
    "387|44461|0|1|53.73455131649651",
    "388|44462|0|1|53.73455131649651",
    "389|44463|0|1|53.73455131649651",
    "390|44464|0|1|53.73455131649651",
    "391|44465|0|1|53.73455131649651",
    "392|44466|0|0|0",
    "393|44467|0|0|0",
    "394|44468|0|0|0",
    "395|44469|0|0|0",
    "396|44470|0|0|0",
    "397|44471|0|0|0",
    "398|44472|0|0|0",
    "399|44473|0|0|0",
    "400|44474|0|0|0",
    "401|44475|0|0|0",
    "402|44476|0|0|0",
    "403|44477|0|0|0",
    "404|44478|0|0|0",
    "405|44479|0|0|0",
    "406|44480|0|0|0",
    "407|44481|0|0|0",
    "408|44482|0|0|0",
    "409|44483|0|0|0",
    "410|44484|0|0|0",
    "411|44485|0|0|0",
    "412|44486|0|0|0",
    "413|44487|0|0|0",
    "414|44488|0|0|0",
    "415|44489|0|0|0",
    "416|44490|0|0|0",
    "417|44491|0|0|0",
    "418|44492|0|0|0",
    "419|44493|0|0|0",
    "420|44494|0|0|0",
    "421|44495|0|0|0",
    "422|44496|0|0|0",
    "423|44497|0|0|0",
    "424|44498|0|0|0",
    "425|44499|0|0|0",
    "426|44500|0|0|0",
    "427|44501|0|0|0",
    "428|44502|0|0|0",
    "429|44503|0|0|0",
    "430|44504|0|0|0",
    "431|44505|0|0|0",
    "432|44506|0|0|0",
    "433|44507|0|0|0",
    "434|44508|0|0|0",
    "435|44509|0|0|0",
    "436|44510|1|1|53.73455131649651",
    "437|44511|0|1|53.73455131649651",
    "438|44512|0|1|53.73455131649651",
    "439|44513|0|1|53.73455131649651",
    "440|44514|0|1|53.73455131649651",
    "441|44515|0|1|53.73455131649651",
    "442|44516|0|1|53.73455131649651",
    "443|44517|0|0|0",
    "444|44518|0|0|0",
    "445|44519|0|0|0",
    "446|44520|0|0|0",
    "447|44521|0|0|0",
    "448|44522|0|0|0",
    "449|44523|0|0|0",
    "450|44524|1|1|53.73455131649651",
    "451|44525|0|1|53.73455131649651",
    "452|44526|1|2|107.469102632993",
    "453|44527|0|2|107.469102632993",
    "454|44528|0|2|107.469102632993",
    "455|44529|0|2|107.469102632993",
    "456|44530|0|2|107.469102632993",
    "457|44531|0|1|53.73455131649651",
    "458|44532|0|1|53.73455131649651",
    "459|44533|0|0|0",
    "460|44534|0|0|0",
    "461|44535|0|0|0",
    "462|44536|0|0|0",
    "463|44537|0|0|0",
    "464|44538|0|0|0"
)

foreach ($line in $rowData) {
    $parts = $line.Split("|")
    $r = [int]$parts[0]
    $a = [double]$parts[1]
    $b = [double]$parts[2]
    $c = [double]$parts[3]
    $d = [double]$parts[4]

    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
}

Write-Host "Done: added rows 386-464"
